# Update "paises" workbook: refresh case counts for a handful of countries,
# swap the Belice/Santa Lucia rows (source list order changed upstream), and
# bump the "last updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp header (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 30 de Mayo de 2020 a las 03:10"

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 1793530
$ws.Range("C4").Value = 25069
$ws.Range("D4").Value = 519569
$ws.Range("E4").Value = 1169419
$ws.Range("G4").Value = 1212
$ws.Range("H4").Value = 104542

# --- Chequia (row 57) ---
$ws.Range("B57").Value = 9196
$ws.Range("C57").Value = 56
$ws.Range("D57").Value = 6500
$ws.Range("E57").Value = 2377

# --- Bahamas (row 173) ---
$ws.Range("B173").Value = 102
$ws.Range("C173").Value = 2
$ws.Range("D173").Value = 48

# --- Belice / Santa Lucia swap places (rows 200 & 201) ---
$ws.Range("A200").Value = "Santa Lucia"
$ws.Range("D200").Value = 18
$ws.Range("H200").Value = 0

$ws.Range("A201").Value = "Belice"
$ws.Range("D201").Value = 16
$ws.Range("H201").Value = 2
